# Append new manifest rows (s103..s119) to the Simulation_Manifest sheet.
# Source data mirrors the XML diff: rows 104-120, columns A-I, with every
# value stored as *text* (the sheet uses inlineStr / text-typed cells
# throughout - even numeric-looking entries like coordinates are text,
# not numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Writes $val into (row, col) as a genuine text cell, even when $val looks
# like an integer (Excel would otherwise auto-convert "937" -> number 937).
# The classic trick is a leading apostrophe (quote-prefix) which forces
# text entry; we then strip the quote-prefix *display* style back to the
# sheet's default "Normal" style so the written cell matches the plain,
# unstyled text cells already present in the sheet.
function Set-TextCell {
    param($Row, $Col, [string]$Val)

    $cell = $ws.Cells.Item($Row, $Col)
    if ($Val -match '^-?[0-9]+$') {
        $cell.Value = "'" + $Val
        $cell.Style = "Normal"
    } else {
        $cell.Value = $Val
    }
}

# row | subject_id | file_name | feedback_1_id | x | y | toleranceA | toleranceB | theta | ratio
$newRows = @'
104|s103|s103_e29_24019-32_1_1.jpeg|meltpatch|937|570|9|9|123|1
105|s104|s104_e30_24019-32_1_2.jpeg|meltpatch|640|1027|9|9|112|1
106|s105|s105_e31_24019-32_1_3.jpeg|meltpatch|201|1333|9|9|53|1
107|s106|s106_e32_24019-32_1_4.jpeg|meltpatch|1806|165|9|9|2|1
108|s107|s107_e33_24019-32_1_5.jpeg|meltpatch|192|405|9|9|142|1
109|s108|s108_e34_24019-32_1_6.jpeg|meltpatch|1688|378|9|9|153|1
110|s109|s109_e35_24019-32_1_7.jpeg|meltpatch|239|1089|9|9|175|1
111|s110|s110_e36_24019-32_1_8.jpeg|meltpatch|1781|1169|9|9|29|1
112|s111|s111_e37_24019-32_2_7.jpeg|meltpatch|1016|616|9|9|91|1
113|s112|s112_e38_24019-32_2_6.jpeg|meltpatch|507|652|9|9|38|1
114|s113|s113_e39_24019-32_2_5.jpeg|meltpatch|1081|1391|9|9|119|1
115|s114|s114_e40_24019-32_2_4.jpeg|meltpatch|1304|311|9|9|18|1
116|s115|s115_e41_24019-32_2_3.jpeg|meltpatch|647|425|9|9|65|1
117|s116|s116_e42_24019-32_2_2.jpeg|meltpatch|1497|1455|9|9|31|1
118|s117|s117_e43_24019-32_2_1.jpeg|meltpatch|1763|182|9|9|150|1
119|s118|s118_e44_24019-32_2_0.jpeg|meltpatch|379|266|9|9|56|1
120|s119|s119_e45_24019-32_3_1.jpeg|meltpatch|61|572|9|9|70|1
'@ -split "`r?`n" | Where-Object { $_.Trim() -ne "" }

foreach ($line in $newRows) {
    $fields = $line -split '\|'
    $row = [int]$fields[0]
    for ($i = 1; $i -lt $fields.Length; $i++) {
        Set-TextCell $row $i $fields[$i]
    }
}
